# Minor update on sequencing assays
#
# 1. dataset_type: rename "Multiplex Ion Beam Imaging" -> "MIBI" and move it
#    up so it sits right before "DESI" (its HRAVS id, HRAVS_0000172, is kept).
# 2. library_preparation_kit: "Custom" kit's NCIt code changes from
#    C126386 -> C65167.
# 3. .metadata: pav:createdOn timestamp bumped.

$wb = $excel.ActiveWorkbook

# --- 1. dataset_type -------------------------------------------------
$ws = $wb.Worksheets.Item("dataset_type")

# Remove the old "Multiplex Ion Beam Imaging" row (row 14) entirely,
# shifting everything below it up by one.
$ws.Rows.Item(14).Delete()

# Insert a fresh row above "DESI" (row 4) for the renamed entry, shifting
# "DESI" and everything after it back down by one.
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Value = "MIBI"
$ws.Cells.Item(4, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000172"

# --- 2. library_preparation_kit --------------------------------------
$ws2 = $wb.Worksheets.Item("library_preparation_kit")
$ws2.Cells.Item(4, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65167"

# --- 3. .metadata ------------------------------------------------------
$ws3 = $wb.Worksheets.Item(".metadata")
$ws3.Cells.Item(2, 3).Value = "2023-10-20T15:01:53-07:00"
